# Split each of the two existing sheets ("sine_ro_ccm", "sine_ro_flash") into
# a FLASH-only sheet and a CCM-only sheet, rename the "FLASH n"/"CCM n" column
# headers down to just "n", and append a new "energy" row with three new
# values to each of the resulting four sheets.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Step 1: duplicate the two original sheets (while they still hold all six
# data columns) so we end up with four sheets in the right order:
#   1: sine_ro_ccm            (-> will become "...code_FLASH")
#   2: sine_ro_ccm (copy)     (-> will become "...code_CCM")
#   3: sine_ro_flash          (-> will become "...code_FLASH")
#   4: sine_ro_flash (copy)   (-> will become "...code_CCM")
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Copy($null, $wb.Worksheets.Item(1))
$wb.Worksheets.Item(3).Copy($null, $wb.Worksheets.Item(3))

# ---------------------------------------------------------------------------
# Step 2: new "energy" data (not present before) to append to each sheet.
# ---------------------------------------------------------------------------
$energy = @{
    1 = @(16535.28305053711, 15757.40281677246, 17004.45761489868)
    2 = @(13688.47531700134, 13365.37822151184, 13196.24352836609)
    3 = @(17069.04348564148, 18640.27112960815, 20366.63442420959)
    4 = @(14420.34118080139, 14619.64906311035, 14933.04797744751)
}

$newNames = @{
    1 = "sine_ro_ccm code_FLASH"
    2 = "sine_ro_ccm code_CCM"
    3 = "sine_ro_flash code_FLASH"
    4 = "sine_ro_flash code_CCM"
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    if ($i % 2 -eq 1) {
        # FLASH sheet: drop the CCM columns (G, E, C), right-to-left so the
        # remaining column letters don't shift out from under us.
        $ws.Columns("G").Delete()
        $ws.Columns("E").Delete()
        $ws.Columns("C").Delete()
    } else {
        # CCM sheet: drop the FLASH columns (F, D, B), right-to-left.
        $ws.Columns("F").Delete()
        $ws.Columns("D").Delete()
        $ws.Columns("B").Delete()
    }

    # Headers become plain "24"/"48"/"72" instead of "FLASH 24"/"CCM 24" etc.
    # Force text storage (they'd otherwise be auto-detected as numbers), then
    # restore the original header formatting (bold/border/center) which the
    # forced text number-format would otherwise disturb.
    $ws.Range("B1:D1").NumberFormat = "@"
    $ws.Range("B1").Value = "24"
    $ws.Range("C1").Value = "48"
    $ws.Range("D1").Value = "72"
    $ws.Range("A2").Copy()
    $ws.Range("B1:D1").PasteSpecial($xlPasteFormats)

    # New row 5: "energy" label + three new numeric values.
    $vals = $energy[$i]
    $ws.Range("A5").Value = "energy"
    $ws.Range("B5").Value = $vals[0]
    $ws.Range("C5").Value = $vals[1]
    $ws.Range("D5").Value = $vals[2]
    $ws.Range("A4").Copy()
    $ws.Range("A5").PasteSpecial($xlPasteFormats)

    $ws.Name = $newNames[$i]
}

$wb.Worksheets.Item(1).Activate()
$excel.CutCopyMode = 0
